$wb = $excel.ActiveWorkbook

# 1. Update the shared "Status" string "Ready for handoff" -> "In Translation"
#    everywhere it appears (Overview, zh-cn, de-de sheets).
foreach ($ws in $wb.Worksheets) {
    $used = $ws.UsedRange
    foreach ($cell in $used.Cells) {
        $text = [string]$cell.Text
        if ($text -eq "Ready for handoff") {
            $cell.Value = "In Translation"
        }
    }
}

# 2. Narrow the "Status" column(s) on each sheet to match the report's
#    new layout (equivalent to a raw stored width of ~13.41 chars).
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Columns.Item(5).ColumnWidth = 12.5
$wsOverview.Columns.Item(6).ColumnWidth = 12.5

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Columns.Item(3).ColumnWidth = 12.5

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Columns.Item(3).ColumnWidth = 12.5
